$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.240.99"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = "'1.648.77"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = "'218.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('E6').Value = '  +1.58%  '

$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('D8').Value = "'0.257"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.84%  '

$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').Value = "'20.24"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.12%  '

$ws.Range('D11').Value = "'0.0848"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '

$ws.Range('D12').Value = "'1.878.19"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.51%  '

$ws.Range('D13').Value = "'1.659.92"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.20%  '

$ws.Range('E14').Value = '  -1.71%  '

$ws.Range('D16').Value = "'67.74"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.91%  '

$ws.Range('D17').Value = "'27.205.03"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '

$ws.Range('D19').Value = "'220.82"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.06%  '

$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('E21').Value = '  +0.33%  '

$ws.Range('E23').Value = '  +1.62%  '

$ws.Range('E24').Value = '  -0.62%  '

$ws.Range('D25').Value = "'148.68"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.61%  '

$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('D27').Value = "'7.41"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('E28').Value = '  +0.17%  '

$ws.Range('D29').Value = "'15.80"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.28%  '

$ws.Range('D30').Value = "'0.0506"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.73%  '

$ws.Range('E31').Value = '  -0.82%  '

$ws.Range('E32').Value = '  -0.54%  '

$ws.Range('D33').Value = "'3.04"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.96%  '

$ws.Range('E34').Value = '  +0.78%  '

$ws.Range('D35').Value = "'1.273.43"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('E37').Value = '  +1.07%  '

$ws.Range('D38').Value = "'0.542"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.55%  '

$ws.Range('D39').Value = "'0.846"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.08%  '

$ws.Range('E40').Value = '  -0.12%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'0.811"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.34%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'2.24"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.79%  '

$ws.Range('E43').Value = '  -0.28%  '

$ws.Range('D44').Value = "'1.788.76"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.58%  '

$ws.Range('D45').Value = "'62.89"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.58%  '

$ws.Range('D46').Value = "'92.35"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.44%  '

$ws.Range('D47').Value = "'1.60"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '

$ws.Range('E48').Value = '  +17.00%  '

$ws.Range('E49').Value = '  -0.72%  '

$ws.Range('D50').Value = "'7.72"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.03%  '

$ws.Range('D51').Value = "'0.0976"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
